$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.001.42"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.909.63"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7824"
$ws.Range("E5").Value = "  +4.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.70"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3155"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.13"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06885"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07960"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7405"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.199"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.79"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "29.991.91"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.90"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.862"
$ws.Range("E18").Value = "  -5.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.21"
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007731"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.849"
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.66"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.258"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1379"
$ws.Range("E27").Value = "  +7.79%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.029"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.370"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.518"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.305"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.077"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05513"
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.253"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7314"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.732"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01929"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.118"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4413"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.77"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8397"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.870"
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.32"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.516"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "980.11"
$ws.Range("E49").Value = "  +8.02%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.14"
$ws.Range("E51").Value = "  -1.27%  "
